# Update Name of Algo
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value  = -21.14920000000001
$ws.Range("A10").Value = -20.51919999999997
$ws.Range("A12").Value = -22.46660000000004
$ws.Range("B13").Value = 5.759699999999999
$ws.Range("A18").Value = -22.31910000000003
